$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the ACC header to include the percent qualifier
$ws.Range("I1").Value = "ACC (%)"

# Row 2 (Slime) - ATK and ACC updated
$ws.Range("F2").Value = 1
$ws.Range("I2").Value = 1.5

# Row 3 (Tanuki) - ATK, DEF and ACC updated
$ws.Range("F3").Value = 0.95
$ws.Range("G3").Value = -0.2
$ws.Range("I3").Value = 1

# Row 4 (The Evil Eye) - ATK, DEF, ACC and EVADE updated
$ws.Range("F4").Value = 1.25
$ws.Range("G4").Value = 0.2
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0

# Update the active selection to F2 (matches the saved sheetView state)
$ws.Range("F2").Select()
